$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$row = $ws.Rows.Item(4)
$row.Insert()
$ws.Range("A4").Value = "Statistics in Practice: Simulation studies as a tool to evaluate and compare the properties of statistical methods – an overview"
$ws.Range("B4").Value = "Willi Sauerbrei"
$ws.Range("A4").WrapText = $true
$ws.Range("B6").Select()
